$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# #167 adding ui ureateam improvement
# Rename "Delhi Daredevils" to "Delhi Capitals"
$ws.Range("A8").Value = "Delhi Capitals"
$ws.Range("A8").Select()
